$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 434, shifting existing rows 434-516 down to 435-517
$ws.Rows.Item(434).Insert()

# Fill in the new row 434 with the new record's data
$ws.Range("A434").Value2 = 5
$ws.Range("B434").Value2 = "Macroferia Regional de Talca"
$ws.Range("C434").Value2 = "Maule"
$ws.Range("D434").Value2 = 44995
$ws.Range("E434").Value2 = 7
$ws.Range("F434").Value2 = 100112032
$ws.Range("G434").Value2 = "Zapallo italiano"
$ws.Range("H434").Value2 = "Sin especificar"
$ws.Range("I434").Value2 = "Primera"
$ws.Range("J434").Value2 = 400
$ws.Range("K434").Value2 = 4000
$ws.Range("L434").Value2 = 4000
$ws.Range("M434").Value2 = 4000
$ws.Range("N434").Value2 = "$/caja 50 unidades"
$ws.Range("O434").Value2 = "Región del Maule"
$ws.Range("P434").Value2 = 80
$ws.Range("Q434").Value2 = 50
$ws.Range("R434").Value2 = "Hortaliza"
